$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.730.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.757.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4469"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3743"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07553"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.127"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.32%  "

$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.213"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.366"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.759.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.20%  "

$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06247"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.184"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5339"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.756.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.316"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.366"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.959.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.223"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09338"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.754"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.644"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02331"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2176"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06151"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6491"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.084"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.199"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.006"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.419"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.28%  "

$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("E46").Value = "  -2.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.753"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.991"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06906"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.24%  "
